$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 12851
$ws1.Range("F5").Value = 83
$ws1.Range("F6").Value = 66
$ws1.Range("F7").Value = 40
$ws1.Range("F8").Value = 20
$ws1.Range("F9").Value = 11
$ws1.Range("F10").Value = 12772
$ws1.Range("F12").Value = 24
$ws1.Range("F13").Value = 8659
$ws1.Range("F14").Value = 7652
$ws1.Range("F16").Value = 90
$ws1.Range("F22").Value = 374
$ws1.Range("F23").Value = 185
$ws1.Range("F24").Value = 15
$ws1.Range("F25").Value = 83
$ws2.Range("F3").Value = 2
$ws4.Range("F4").Value = 12851
$ws4.Range("F6").Value = 83
$ws4.Range("F7").Value = 66
$ws4.Range("F8").Value = 40
$ws4.Range("F9").Value = 20
$ws4.Range("F10").Value = 11
$ws4.Range("F11").Value = 12772
$ws4.Range("F13").Value = 24
$ws4.Range("F14").Value = 8659
$ws4.Range("F15").Value = 7652
$ws4.Range("F17").Value = 90
$ws4.Range("F23").Value = 2
$ws4.Range("F24").Value = 374
$ws4.Range("F25").Value = 185
$ws4.Range("F26").Value = 15
$ws4.Range("F27").Value = 83
